$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder / rewrite the header row (row 1) -----------------------------
# Target header layout (columns A..O). Column K duplicates "S. No." the same
# way the source workbook did after the author's edit.
$nbsp = [char]0xA0
$headers = @(
    "S. No.",
    "INSTITUTE CODE",
    "INSTITUTE NAME ",
    ("MOBILE" + $nbsp),
    "ALTERNATIVE MOBILE NUMBER",
    "EMAIL-ID",
    "WEBSITE",
    "ADDRESS",
    ("DISTRICT" + $nbsp),
    ("PINCODE" + $nbsp),
    "S. No.",
    "BRANCH CODE",
    "BRANCH NAME",
    "TOTAL SEATS",
    "INSTITUTE TYPE"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Normalise formatting across the whole header row (A1 already carries the
# correct look from the original template: bordered, shaded, centred,
# wrapped Garamond cell) - copy it across so every header cell (including
# the brand new O1) shares the exact same style, instead of O1 falling back
# to the row's bare default style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1:O1").PasteSpecial(-4122) | Out-Null

# --- New sample data rows (2-5): columns A and B only ----------------------
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = 2
}

# --- Header row formatting --------------------------------------------------
$ws.Range("A1:O1").Font.Bold = $true
$ws.Rows(1).RowHeight = 31

# --- Column widths: auto-fit to content (closest achievable match to the
#     bestFit widths Excel itself would compute) -----------------------------
$ws.Columns("A:O").AutoFit() | Out-Null

# --- View / selection -------------------------------------------------------
$ws.Range("C1").Select() | Out-Null
$ws.Range("H10").Select() | Out-Null

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1
